$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = -7.711
$ws.Range("B3").Value = 6.265
$ws.Range("D3").Value = -7.388
$ws.Range("B4").Value = 7.093000000000001
$ws.Range("D9").Value = -7.028
$ws.Range("A11").Value = -21.452
$ws.Range("A12").Value = -21.702
$ws.Range("B14").Value = 6.066
$ws.Range("A15").Value = -21.209
$ws.Range("D15").Value = -8.111000000000001
$ws.Range("D19").Value = -8.159000000000001
$ws.Range("D20").Value = -7.765000000000001
$ws.Range("D25").Value = -7.683
$ws.Range("B26").Value = 6.225
$ws.Range("A27").Value = -20.919
$ws.Range("D27").Value = -8.231
$ws.Range("A28").Value = -21.527
$ws.Range("D28").Value = -8.004999999999999
$ws.Range("D30").Value = -7.075999999999999
$ws.Range("A31").Value = -21.374
$ws.Range("B31").Value = 6.108999999999999
$ws.Range("A32").Value = -20.963
$ws.Range("D32").Value = -8.149000000000001
$ws.Range("B35").Value = 7.216000000000001
$ws.Range("A36").Value = -20.925
$ws.Range("B37").Value = 7.065
$ws.Range("A38").Value = -19.741
$ws.Range("B39").Value = 7.615
$ws.Range("B40").Value = 8.211
$ws.Range("D44").Value = -8.242000000000001
$ws.Range("B45").Value = 5.87
$ws.Range("A46").Value = -21.266
$ws.Range("D47").Value = -7.489
$ws.Range("B52").Value = 5.518
$ws.Range("A54").Value = -21.856
$ws.Range("A55").Value = -22.21
$ws.Range("A56").Value = -21.803
$ws.Range("B57").Value = 5.332000000000001
$ws.Range("D58").Value = -8.146000000000001
$ws.Range("D62").Value = -8.087
$ws.Range("A67").Value = -21.586
$ws.Range("A69").Value = -21.637
$ws.Range("A72").Value = -21.567
$ws.Range("A73").Value = -20.621
$ws.Range("D77").Value = -7.73
$ws.Range("D78").Value = -7.812
$ws.Range("B81").Value = 6.392000000000001
$ws.Range("A83").Value = -20.531
$ws.Range("B83").Value = 6.207999999999999
$ws.Range("D84").Value = -8.148
$ws.Range("A86").Value = -21.911
$ws.Range("D89").Value = -6.790999999999999
$ws.Range("A91").Value = -21.522
$ws.Range("D91").Value = -6.929
$ws.Range("D92").Value = -6.741
$ws.Range("A93").Value = -21.665
$ws.Range("D96").Value = -7.473000000000001
$ws.Range("A99").Value = -20.547
$ws.Range("B100").Value = 5.558
$ws.Range("B102").Value = 7.499000000000001
$ws.Range("D102").Value = -7.986999999999999
